# #5: insurance, claim, debt, investment done
#
# The "保險" (insurance) sheet had its header row (row 1) mistakenly
# filled with copies of the first data row instead of the field names,
# and it was missing the usual metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index)
# that every other property-table sheet in this workbook carries.
# Also relabels the stray "otherbonds" category value (used on the
# "具有相當價值之財產" sheet) to the correct "antique" label.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 具有相當價值之財產 (property of considerable value): otherbonds -> antique
# ---------------------------------------------------------------
$wsValuable = $wb.Worksheets.Item("具有相當價值之財產")
$wsValuable.Range("F2").Value = "antique"

# ---------------------------------------------------------------
# 保險 (insurance) sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("保險")

# --- fix header row: it held stray data instead of field labels ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"

# --- add the standard metadata header columns (E1:K1) -------------
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# header formatting to match the rest of the header row (bold, boxed,
# centered) instead of the plain default style the new cells start with
$headerRange = $ws.Range("E1:K1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- row 2 (index 91): add property metadata columns E2:K2 --------
# format the "date" column as text first so the ISO-looking
# "2011-12-22" string isn't auto-converted into a date serial number
$ws.Range("G2:G3").NumberFormat = "@"

$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("G2").Value = "2011-12-22"
$ws.Range("H2").Value = "徐耀昌"
$ws.Range("I2").Value = 921
$ws.Range("J2").Value = "tmpd3a41"
$ws.Range("K2").Value = 91

# --- row 3 (index 92): add property metadata columns E3:K3 --------
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("G3").Value = "2011-12-22"
$ws.Range("H3").Value = "徐耀昌"
$ws.Range("I3").Value = 921
$ws.Range("J3").Value = "tmpd3a41"
$ws.Range("K3").Value = 92
